$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing readings for 2023-05-30 (row 6 and row 12) ---
# C6: shift-1 time reading corrected to 16:15 (0.67708333333333337)
$ws.Range("C6").Value = 0.67708333333333337
# C12: shift-2 time reading corrected to 23:00 (0.95833333333333337)
$ws.Range("C12").Value = 0.95833333333333337

# --- Append the 2023-05-31 readings (rows 13-19) ---
# Seed the date/time formatting by copying the formats already used
# for the 2023-05-30 block, then overwrite the values.

# Row 13
$ws.Range("A2").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 45077
$ws.Range("B13").Value = 1
$ws.Range("C2").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = 0.40625
$ws.Range("D13").Value = 1130

# Row 14
$ws.Range("A2").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 45077
$ws.Range("B14").Value = 1
$ws.Range("C2").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = 0.45833333333333331
$ws.Range("D14").Value = 1140

# Row 15
$ws.Range("A2").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = 45077
$ws.Range("B15").Value = 1
$ws.Range("C2").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 0.52083333333333337
$ws.Range("D15").Value = 1340

# Row 16
$ws.Range("A2").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 45077
$ws.Range("B16").Value = 1
$ws.Range("C2").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 0.58333333333333337
$ws.Range("D16").Value = 1232

# Row 17
$ws.Range("A2").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 45077
$ws.Range("B17").Value = 1
$ws.Range("C2").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 0.67708333333333337
$ws.Range("D17").Value = 1342

# Row 18 (date cell additionally gets an explicit black font)
$ws.Range("A2").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = 45077
$ws.Range("A18").Font.Color = 0
$ws.Range("B18").Value = 2
$ws.Range("C2").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 0.71527777777777779
$ws.Range("D18").Value = 1156

# Row 19 (date cell additionally gets an explicit black font)
$ws.Range("A2").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = 45077
$ws.Range("A19").Font.Color = 0
$ws.Range("B19").Value = 2
$ws.Range("C2").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C19").Value = 0.77083333333333337
$ws.Range("D19").Value = 1116

# --- Move the active selection to match where the user continued entry ---
$ws.Range("E12").Select() | Out-Null
